$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 126
$ws.Range("I11").Value = 126
$ws.Range("K11").Value = 126
$ws.Range("M11").Value = 14
$ws.Range("H39").Value = 402.77777
$ws.Range("I39").Value = 237.5
$ws.Range("J39").Value = 733.3333
$ws.Range("K39").Value = 712.5
$ws.Range("L39").Value = 2199.9999
$ws.Range("M39").Value = -416.5
$ws.Range("N39").Value = -2791.9999
$ws.Range("H41").Value = 316.5
$ws.Range("I41").Value = 316.5
$ws.Range("K41").Value = 316.5
$ws.Range("M41").Value = 123.5
$ws.Range("H99").Value = 400
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H132").Value = 1029.75
$ws.Range("I132").Value = 987.0454999999999
$ws.Range("J132").Value = 1499.5
$ws.Range("K132").Value = 2961.1365
$ws.Range("L132").Value = 4498.5
$ws.Range("M132").Value = -431.1364999999996
$ws.Range("N132").Value = -9558.5
$ws.Range("H137").Value = 1524.2667
$ws.Range("I137").Value = 1096.2222
$ws.Range("K137").Value = 3288.6666
$ws.Range("M137").Value = -738.6665999999996

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 18500
$ws.Range("J23").Value = 18500
$ws.Range("L23").Value = 18500
$ws.Range("N23").Value = -19018
$ws.Range("H27").Value = 17500
$ws.Range("J27").Value = 17500
$ws.Range("L27").Value = 17500
$ws.Range("N27").Value = -17868
$ws.Range("H32").Value = 6747.884
$ws.Range("I32").Value = 4604.4326
$ws.Range("K32").Value = 4604.4326
$ws.Range("M32").Value = -4317.4326
$ws.Range("H61").Value = 1384.8182
$ws.Range("I61").Value = 1433.3
$ws.Range("J61").Value = 900
$ws.Range("K61").Value = 1433.3
$ws.Range("L61").Value = 900
$ws.Range("M61").Value = -1221.3
$ws.Range("N61").Value = -1324
$ws.Range("H74").Value = 1221.5385
$ws.Range("I74").Value = 1023
$ws.Range("K74").Value = 1023
$ws.Range("M74").Value = -149
$ws.Range("H77").Value = 1221.5385
$ws.Range("I77").Value = 1023
$ws.Range("K77").Value = 5115
$ws.Range("M77").Value = -747
$ws.Range("H110").Value = 107983.164
$ws.Range("I110").Value = 157977.25
$ws.Range("K110").Value = 157977.25
$ws.Range("M110").Value = -155932.25
$ws.Range("H122").Value = 2597.8125
$ws.Range("I122").Value = 2361.7856
$ws.Range("J122").Value = 4250
$ws.Range("K122").Value = 7085.3568
$ws.Range("L122").Value = 12750
$ws.Range("M122").Value = -4635.3568
$ws.Range("N122").Value = -17650
$ws.Range("H132").Value = 1056.3334
$ws.Range("I132").Value = 646
$ws.Range("J132").Value = 1466.6666
$ws.Range("K132").Value = 1938
$ws.Range("L132").Value = 4399.9998
$ws.Range("M132").Value = 592
$ws.Range("N132").Value = -9459.9998
$ws.Range("H136").Value = 1384.8182
$ws.Range("I136").Value = 1433.3
$ws.Range("J136").Value = 900
$ws.Range("K136").Value = 4299.9
$ws.Range("L136").Value = 2700
$ws.Range("M136").Value = -1749.9
$ws.Range("N136").Value = -7800

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 2833.3333
$ws.Range("I30").Value = 2833.3333
$ws.Range("K30").Value = 2833.3333
$ws.Range("M30").Value = -2708.3333
$ws.Range("H99").Value = 2853.2
$ws.Range("I99").Value = 3299.5
$ws.Range("K99").Value = 3299.5
$ws.Range("M99").Value = -1801.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3490.6
$ws.Range("I16").Value = 3426
$ws.Range("J16").Value = 3749
$ws.Range("K16").Value = 3426
$ws.Range("L16").Value = 3749
$ws.Range("M16").Value = -3139
$ws.Range("N16").Value = -4323
$ws.Range("H31").Value = 3030.742
$ws.Range("I31").Value = 2526.0386
$ws.Range("J31").Value = 5655.2
$ws.Range("K31").Value = 2526.0386
$ws.Range("L31").Value = 5655.2
$ws.Range("M31").Value = -2231.0386
$ws.Range("N31").Value = -6245.2
$ws.Range("H34").Value = 3030.742
$ws.Range("I34").Value = 2526.0386
$ws.Range("J34").Value = 5655.2
$ws.Range("K34").Value = 2526.0386
$ws.Range("L34").Value = 5655.2
$ws.Range("M34").Value = -2324.0386
$ws.Range("N34").Value = -6059.2
$ws.Range("H58").Value = 1812.7693
$ws.Range("I58").Value = 1119.091
$ws.Range("K58").Value = 1119.091
$ws.Range("M58").Value = -916.0909999999999
$ws.Range("H99").Value = 14467.19
$ws.Range("I99").Value = 11566.375
$ws.Range("J99").Value = 16252.308
$ws.Range("K99").Value = 11566.375
$ws.Range("L99").Value = 16252.308
$ws.Range("M99").Value = -10068.375
$ws.Range("N99").Value = -19248.308
$ws.Range("H113").Value = 3490.6
$ws.Range("I113").Value = 3426
$ws.Range("J113").Value = 3749
$ws.Range("K113").Value = 3426
$ws.Range("L113").Value = 3749
$ws.Range("M113").Value = -1256
$ws.Range("N113").Value = -8089
$ws.Range("H126").Value = 14467.19
$ws.Range("I126").Value = 11566.375
$ws.Range("J126").Value = 16252.308
$ws.Range("K126").Value = 34699.125
$ws.Range("L126").Value = 48756.924
$ws.Range("M126").Value = -32229.125
$ws.Range("N126").Value = -53696.924
$ws.Range("H132").Value = 1788.5
$ws.Range("I132").Value = 1616.25
$ws.Range("K132").Value = 4848.75
$ws.Range("M132").Value = -2318.75
$ws.Range("H134").Value = 1868.6945
$ws.Range("I134").Value = 1685.9333
$ws.Range("K134").Value = 5057.7999
$ws.Range("M134").Value = -2522.7999
$ws.Range("H136").Value = 1812.7693
$ws.Range("I136").Value = 1119.091
$ws.Range("K136").Value = 3357.273
$ws.Range("M136").Value = -807.2729999999997

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 188.33333
$ws.Range("I23").Value = 188.33333
$ws.Range("K23").Value = 564.99999
$ws.Range("M23").Value = -329.99999
$ws.Range("H39").Value = 2633.3333
$ws.Range("J39").Value = 2950
$ws.Range("L39").Value = 8850
$ws.Range("N39").Value = -9438
$ws.Range("H97").Value = 8371.25
$ws.Range("I97").Value = 7995
$ws.Range("K97").Value = 23985
$ws.Range("M97").Value = -23489
$ws.Range("H112").Value = 4990
$ws.Range("J112").Value = 4990
$ws.Range("L112").Value = 14970
$ws.Range("N112").Value = -17186
$ws.Range("H130").Value = 3581.1667
$ws.Range("I130").Value = 999
$ws.Range("J130").Value = 4872.25
$ws.Range("K130").Value = 2997
$ws.Range("L130").Value = 14616.75
$ws.Range("M130").Value = 2023
$ws.Range("N130").Value = -24656.75
$ws.Range("H131").Value = 1122.909
$ws.Range("I131").Value = 829.3333
$ws.Range("J131").Value = 1152.2667
$ws.Range("K131").Value = 2487.9999
$ws.Range("L131").Value = 3456.800099999999
$ws.Range("M131").Value = 2552.0001
$ws.Range("N131").Value = -13536.8001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2494.5715
$ws.Range("I102").Value = 2494.5715
$ws.Range("K102").Value = 2494.5715
$ws.Range("M102").Value = -872.5715
$ws.Range("H132").Value = 2412.875
$ws.Range("I132").Value = 2360.375
$ws.Range("J132").Value = 2517.875
$ws.Range("K132").Value = 7081.125
$ws.Range("L132").Value = 7553.625
$ws.Range("M132").Value = -4551.125
$ws.Range("N132").Value = -12613.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2589.2727
$ws.Range("I40").Value = 2649.3
$ws.Range("K40").Value = 2649.3
$ws.Range("M40").Value = -2513.3
$ws.Range("H61").Value = 3466.6155
$ws.Range("I61").Value = 3338.8333
$ws.Range("K61").Value = 3338.8333
$ws.Range("M61").Value = -3136.8333
$ws.Range("H82").Value = 2200.4
$ws.Range("J82").Value = 4001.5
$ws.Range("L82").Value = 4001.5
$ws.Range("N82").Value = -4723.5
$ws.Range("H85").Value = 2200.4
$ws.Range("J85").Value = 4001.5
$ws.Range("L85").Value = 4001.5
$ws.Range("N85").Value = -6497.5
$ws.Range("H113").Value = 3466.6155
$ws.Range("I113").Value = 3338.8333
$ws.Range("K113").Value = 3338.8333
$ws.Range("M113").Value = -1168.8333
$ws.Range("H115").Value = 70000
$ws.Range("J115").Value = 70000
$ws.Range("L115").Value = 70000
$ws.Range("N115").Value = -72350
$ws.Range("H132").Value = 5601
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 10005
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 30015
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -35075
$ws.Range("H136").Value = 6172.875
$ws.Range("I136").Value = 4399.5
$ws.Range("J136").Value = 7946.25
$ws.Range("K136").Value = 13198.5
$ws.Range("L136").Value = 23838.75
$ws.Range("M136").Value = -10648.5
$ws.Range("N136").Value = -28938.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2059.8
$ws.Range("I96").Value = 1699.75
$ws.Range("J96").Value = 3500
$ws.Range("K96").Value = 1699.75
$ws.Range("L96").Value = 3500
$ws.Range("M96").Value = -326.75
$ws.Range("N96").Value = -6246
$ws.Range("H132").Value = 7181.0713
$ws.Range("I132").Value = 4176.857
$ws.Range("J132").Value = 10185.286
$ws.Range("K132").Value = 12530.571
$ws.Range("L132").Value = 30555.858
$ws.Range("M132").Value = -10000.571
$ws.Range("N132").Value = -35615.858
$ws.Range("H136").Value = 1538.9166
$ws.Range("I136").Value = 1176.5264
$ws.Range("J136").Value = 2916
$ws.Range("K136").Value = 3529.5792
$ws.Range("L136").Value = 8748
$ws.Range("M136").Value = -979.5792000000001
$ws.Range("N136").Value = -13848
